$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): insert "gender" and "status" after "name", shift the
# original email..levelName set right, replace the old pfaName /
# pensionAccountNumber columns with the new BPJS columns, and append npwp +
# levelName at the end. Written as literal values (not a column Insert) so
# row 2's existing B2 cell is left exactly where it is. ---
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "gender"
$ws.Range("C1").Value = "status"
$ws.Range("D1").Value = "email"
$ws.Range("E1").Value = "designation"
$ws.Range("F1").Value = "department"
$ws.Range("G1").Value = "stateResidence"
$ws.Range("H1").Value = "bankName"
$ws.Range("I1").Value = "accountNumber"
$ws.Range("J1").Value = "bpjsKetenagakerjaanNumber"
$ws.Range("K1").Value = "bpjsKesehatanNumber"
$ws.Range("L1").Value = "npwp"
$ws.Range("M1").Value = "levelName"

# --- Column widths ---
# Existing column J grows to fit the longer header text; K/L/M are brand new
# columns that need an explicit width. (ColumnWidth is in "characters"; the
# engine adds a fixed 5/6-character padding on top when it stores the OOXML
# <col width>, so subtract that back out to land on the authored widths.)
$pad = 5 / 6
$ws.Columns(10).ColumnWidth = 30.85546875 - $pad
$ws.Columns(11).ColumnWidth = 27 - $pad
$ws.Columns(12).ColumnWidth = 17.28515625 - $pad
$ws.Columns(13).ColumnWidth = 17.7109375 - $pad

# --- Selection moves to P6 ---
[void]$ws.Range("P6").Select()
